$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.539005637168884
$ws.Range("B1").Value = 2.381497144699097
$ws.Range("C1").Value = 4.335558891296387
$ws.Range("D1").Value = 1.828066468238831
$ws.Range("E1").Value = 0.8140384554862976
